$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "D", "E", "F", "G", "H", "I", "K", "N")

$rowData = @(
    @(0.07838067010828809, 0.01115829819396552, 0.4197278238190023, 0.774093437100376, 0.6308869797421579, 0.6806869475433075, 0.4094081953048612, 0.9677831610811154, 1.169854919506722),
    @(0.06969514161845325, 0.01105465081561974, 0.3660050950540494, 0.7541389606044646, 0.6101860811634907, 0.6765671856120861, 0.4149996401624207, 0.845085151049517, 1.186028645880063),
    @(0.06436407454299342, 0.01099207512753608, 0.3331304009831655, 0.742510040713924, 0.5980410425376732, 0.6744956608340829, 0.4186535000866929, 0.7697442166640656, 1.196451072706473),
    @(0.06219229207425769, 0.010966845912181, 0.3197593613606529, 0.7379269206433037, 0.5932328352636773, 0.6737662527809789, 0.4201979122179598, 0.7390403713625062, 1.200821805120041),
    @(0.06183171557375999, 0.0109626730851371, 0.317540593218979, 0.7371752768468411, 0.5924429175202448, 0.6736520538913311, 0.420457707040832, 0.7339418885841553, 1.201555017761328),
    @(0.06433478214719912, 0.01099173377572171, 0.3329499735768593, 0.742447601754165, 0.5979756280528505, 0.6744853597047467, 0.4186741041978657, 0.7693301422466732, 1.19650951803128),
    @(0.07538562454283237, 0.0111223414932784, 0.4011797161705886, 0.7670832032334971, 0.6236312420156906, 0.6791711335311845, 0.4112902857451068, 0.9254771453293245, 1.175329491342513),
    @(0.0970636411496173, 0.01138680115314017, 0.5359720803851928, 0.8203851613439639, 0.6784847933926414, 0.6920146808389802, 0.3985642423974021, 1.231703733799122, 1.137702168370005),
    @(0.1129863452014916, 0.01158609776535258, 0.6357734190404898, 0.8626593198027592, 0.7216390780214681, 0.7037096497254538, 0.3902867102690024, 1.456799196587042, 1.112445342014261),
    @(0.1202273777192033, 0.01167784171142472, 0.6813788983944846, 0.8825825391413247, 0.7419092097774751, 0.7095273497426717, 0.3867545862687116, 1.559247592317831, 1.101475301068469),
    @(0.1229688618867755, 0.01171273810550844, 0.6986808309565049, 0.8902276992125735, 0.7496782937488717, 0.7118024437838244, 0.3854506886006401, 1.598050889406409, 1.097396084045686),
    @(0.1223784614525982, 0.01170521565729388, 0.6949530795721159, 0.8885766824193269, 0.7480009149234377, 0.7113092487163044, 0.385730008767478, 1.58969352026412, 1.098271281875425),
    @(0.1204529332421771, 0.01168070955637646, 0.6828016799019281, 0.8832094872128948, 0.7425465009384027, 0.7097130760561186, 0.3866466391279442, 1.562439790888732, 1.101138199123556),
    @(0.1192734154028159, 0.0116657190260625, 0.6753628584205416, 0.8799350678424389, 0.7392176953280227, 0.7087447732544376, 0.3872124852252554, 1.545747205400517, 1.102904029094316),
    @(0.1125130628157223, 0.01158012379554307, 0.6327973471332484, 0.8613713181524503, 0.7203273498212468, 0.7033395051656157, 0.3905222459113773, 1.450105064993863, 1.113172729627004),
    @(0.1083650710653075, 0.01152789054872727, 0.6067392166903716, 0.8501611943803198, 0.7089033500712105, 0.7001513511777659, 0.3926125067054862, 1.391445512564246, 1.119605426006656),
    @(0.1059790532400058, 0.01149794931009396, 0.5917704401600616, 0.8437785424749507, 0.7023926356823722, 0.6983644152891202, 0.3938367300199221, 1.357710864594935, 1.123354219968475),
    @(0.1051711601778038, 0.01148782926152947, 0.5867054836943879, 0.8416286367291832, 0.7001984960170944, 0.6977674132540699, 0.3942550002254279, 1.34628968467365, 1.124631885627169),
    @(0.1088066543617288, 0.01153344032448445, 0.6095111462950626, 0.8513477842244015, 0.7101132279273656, 0.7004858877733966, 0.3923877214445213, 1.397689421179905, 1.118915594430169),
    @(0.1210185238073649, 0.0116879033899977, 0.6863699477223122, 0.8847832216646481, 0.7441460537494322, 0.7101799509980538, 0.3863764888975219, 1.570444638681522, 1.100294081599865),
    @(0.1289964559959316, 0.01178975711377461, 0.7367901127820176, 0.9072225052935892, 0.7669324881057378, 0.7169358255277984, 0.3826439256791332, 1.683399066395509, 1.088560517839666),
    @(0.1247388478524982, 0.01173531334669065, 0.7098618061571216, 0.8951921366799809, 0.754720717683, 0.7132914676804774, 0.3846180918368649, 1.623108372960587, 1.094782916065992),
    @(0.1086070185311883, 0.01153093099518898, 0.6082579185514021, 0.8508111332970714, 0.7095660644030772, 0.7003345005998938, 0.3924892767672823, 1.394866582849886, 1.119227309834435),
    @(0.09119927091094837, 0.01131437575318373, 0.4993839866810816, 0.8054239697196977, 0.6631504762681004, 0.6881457377076856, 0.4018189057099235, 1.1488489659979, 1.147462484397048)
)

$startRow = 2
for ($i = 0; $i -lt $rowData.Length; $i++) {
    $rowVals = $rowData[$i]
    $rowNum = $startRow + $i
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $col = $cols[$j]
        $addr = "$col$rowNum"
        $ws.Range($addr).Value = $rowVals[$j]
    }
}

Write-Host "Done updating pl_mw values for 380 kV case."